$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSub = "dataset1\aa\opt"
$newSub = "dataset1\.tar\opt"

foreach ($r in 2..5) {
    foreach ($col in @("A", "B")) {
        $cell = $ws.Range("$col$r")
        $text = $cell.Value()
        $text = $text.Replace($oldSub, $newSub)
        $cell.Value = $text
    }
}
